$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture current (before) values for rows 3, 4, 5 in the columns that change
$r3D = $ws.Range("D3").Value2
$r3I = $ws.Range("I3").Value2
$r3J = $ws.Range("J3").Value2
$r3K = $ws.Range("K3").Value2
$r3L = $ws.Range("L3").Value2
$r3M = $ws.Range("M3").Value2
$r3P = $ws.Range("P3").Value2

$r4D = $ws.Range("D4").Value2
$r4I = $ws.Range("I4").Value2
$r4J = $ws.Range("J4").Value2
$r4K = $ws.Range("K4").Value2
$r4L = $ws.Range("L4").Value2
$r4M = $ws.Range("M4").Value2
$r4P = $ws.Range("P4").Value2

$r5D = $ws.Range("D5").Value2
$r5I = $ws.Range("I5").Value2
$r5J = $ws.Range("J5").Value2
$r5K = $ws.Range("K5").Value2
$r5L = $ws.Range("L5").Value2
$r5M = $ws.Range("M5").Value2
$r5P = $ws.Range("P5").Value2

# Cyclic rotation: row3 -> row4, row4 -> row5, row5 -> row3
$ws.Range("D4").Value2 = $r3D
$ws.Range("I4").Value2 = $r3I
$ws.Range("J4").Value2 = $r3J
$ws.Range("K4").Value2 = $r3K
$ws.Range("L4").Value2 = $r3L
$ws.Range("M4").Value2 = $r3M
$ws.Range("P4").Value2 = $r3P

$ws.Range("D5").Value2 = $r4D
$ws.Range("I5").Value2 = $r4I
$ws.Range("J5").Value2 = $r4J
$ws.Range("K5").Value2 = $r4K
$ws.Range("L5").Value2 = $r4L
$ws.Range("M5").Value2 = $r4M
$ws.Range("P5").Value2 = $r4P

$ws.Range("D3").Value2 = $r5D
$ws.Range("I3").Value2 = $r5I
$ws.Range("J3").Value2 = $r5J
$ws.Range("K3").Value2 = $r5K
$ws.Range("L3").Value2 = $r5L
$ws.Range("M3").Value2 = $r5M
$ws.Range("P3").Value2 = $r5P
